$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of column J (|S*|/n)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Summary labels (column A) and formulas (column B)
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Style the first summary cell (bold, 12pt, vertically centered) then
# propagate the same formatting to the remaining summary cells.
$c = $ws.Range("B14")
$c.Font.Bold = $true
$c.Font.Size = 12
$c.VerticalAlignment = -4108

$c.Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row heights for the summary rows
$ws.Rows(14).RowHeight = 15.6
$ws.Rows(15).RowHeight = 15.6
$ws.Rows(16).RowHeight = 15.6
$ws.Rows(17).RowHeight = 15.6

# Selection matching the author's last saved state
$ws.Range("A14:B17").Select() | Out-Null

# Print setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
